# isim tamam, eğitmen tamam değil
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second person (row 2): mustafa aksfmasl <mr.turran@gmail.com>
$ws.Range("A2").Value = "mustafa"
$ws.Range("B2").Value = "aksfmasl"
$ws.Range("C2").Value = "mr.turran@gmail.com"

# Third person (row 3): Umut güzel <mr.turran@gmail.com>
$ws.Range("A3").Value = "Umut"
$ws.Range("B3").Value = "güzel"
$ws.Range("C3").Value = "mr.turran@gmail.com"

# Turn the e-mail addresses into mailto: hyperlinks (row1 header included)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:mr.turran@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:onurturan.t@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:mr.turran@gmail.com")

# Match the row height that Excel recalculates for the header row once its
# font changes because of the new hyperlink style
$ws.Rows.Item(1).RowHeight = 15

# Leave the selection where the author left it before saving
[void]$ws.Range("F11").Select()
